# Fikser avrunding av dagsats for 6G, og gradert dagsats av 6G
#
# The "6G daily rate" (C14) and the resulting "daily rate capped at 6G"
# (C15) on the "Multiple Employers" sheet were not rounded, which produced
# a non-integer daily rate. Wrap both formulas in ROUND(..., 0) so the
# amounts are rounded to whole kroner. Every other cell on the sheet that
# changed value in the diff (C16, D17:G19, F22:G22, D23:E23, C28, D33,
# C34, C36, ...) is a downstream formula that recalculates automatically
# once these two inputs change - nothing else needs to be touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Multiple Employers")

$ws.Range("C14").Formula = "=ROUND(6*C13/260, 0)"
$ws.Range("C15").Formula = "=ROUND(C12*C14,0)"

# Match the author's final cell selection on this sheet.
$ws.Range("C16").Select()
